$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 376, pushing existing rows 376-398 down to 377-399.
$ws.Range("A376").EntireRow.Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Cells.Item(376, 1).Value = 5
$ws.Cells.Item(376, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(376, 3).Value = "Maule"
$ws.Cells.Item(376, 4).Value = 45013
$ws.Cells.Item(376, 5).Value = 7
$ws.Cells.Item(376, 6).Value = 100112045
$ws.Cells.Item(376, 7).Value = "Zapallo"
$ws.Cells.Item(376, 8).Value = "Camote"
$ws.Cells.Item(376, 9).Value = "1a (cosecha)"
$ws.Cells.Item(376, 10).Value = 900
$ws.Cells.Item(376, 11).Value = 250
$ws.Cells.Item(376, 12).Value = 250
$ws.Cells.Item(376, 13).Value = 250
$ws.Cells.Item(376, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(376, 15).Value = "Región del Maule"
$ws.Cells.Item(376, 16).Value = 250
$ws.Cells.Item(376, 17).Value = 1
$ws.Cells.Item(376, 18).Value = "Hortaliza"
